$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 3-7 (columns A, B, E, F, G, H, Q, R)
$data = @{
    3 = @{ A = 79930402; B = 89557; E = 1588; F = "Violmussling"; G = "Trichaptum laricinum"; H = "(P.Karst.) Ryvarden"; Q = 400486.0196595828; R = 6751207.17080476 }
    4 = @{ A = 79930409; B = 77506; E = 6425; F = "Garnlav"; G = "Alectoria sarmentosa"; H = "(Ach.) Ach."; Q = 400637.2251421487; R = 6751170.796163691 }
    5 = @{ A = 79930410; B = 85703; E = 510; F = "Doftskinn"; G = "Cystostereum murrayi"; H = "(Berk. & M.A. Curtis.) Pouzar"; Q = 400228.82899852; R = 6751134.972044618 }
    6 = @{ A = 79930407; B = 76504; E = 314; F = "Vitskaftad svartspik"; G = "Chaenothecopsis viridialba"; H = "(Kremp.) A.F.W.Schmidt"; Q = 400891.1128622342; R = 6751229.963137357 }
    7 = @{ A = 79930405; B = 73693; E = 6440; F = "Vitgrynig nållav"; G = "Chaenotheca subroscida"; H = "(Eitner) Zahlbr."; Q = 400891.1128622342; R = 6751229.963137357 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
